# Teaching-experience workbook update: add a new "PSY 3010: Psycological
# Statistics" undergrad course row, tidy up a couple of date labels, and
# drop a duplicated description line in the STAT 1040 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shorten the EDUC 6600 "when" label.
$ws.Range("C7").Value = "Spring/Sum 2015-present"

# 2. Tidy up the MATH 1050 and MATH 1010 "when" labels.
$ws.Range("C20").Value = "Fall 2004"
$ws.Range("C24").Value = "Summer 2006"

# 3. Insert two new rows before the STAT 2000 block (old row 13) for the new
#    PSY 3010 course entry.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "undergrad"
$ws.Range("B13").Value = "Psychology Department"
$ws.Range("C13").Value = "Spring 2022"
$ws.Range("D13").Value = "PSY 3010: Psycological Statistics"
$ws.Range("E13").Value = "Utah State University"
$ws.Range("F13").Value = "Descriptive and inferential statistical methods."
$ws.Range("F14").Value = "Focus on behavior statistics and statistical applications."

# 4. Remove the duplicated "Descriptive and inferential statistical methods."
#    line in the STAT 1040 block (now at row 18, just below the STAT 1040
#    header row 17, after the insert above shifted everything down by 2).
$ws.Rows.Item(18).Delete()

# 5. Update the active selection to match the saved view.
$ws.Range("B7").Select()
